$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 254-265: update changed values (date, volume/price stats, origin) per diff
$ws.Range("D254").Value = 44509
$ws.Range("J254").Value = 200
$ws.Range("K254").Value = 900
$ws.Range("L254").Value = 900
$ws.Range("M254").Value = 900
$ws.Range("O254").Value = 'Región Metropolitana'
$ws.Range("P254").Value = 900

$ws.Range("D255").Value = 44509
$ws.Range("J255").Value = 800
$ws.Range("K255").Value = 800
$ws.Range("L255").Value = 900
$ws.Range("M255").Value = 850
$ws.Range("O255").Value = 'Región del Maule'
$ws.Range("P255").Value = 850

$ws.Range("D256").Value = 44421
$ws.Range("J256").Value = 4300
$ws.Range("K256").Value = 800
$ws.Range("L256").Value = 1000
$ws.Range("M256").Value = 884
$ws.Range("O256").Value = 'Región Metropolitana'
$ws.Range("P256").Value = 884

$ws.Range("D257").Value = 44421
$ws.Range("J257").Value = 1450
$ws.Range("K257").Value = 800
$ws.Range("L257").Value = 900
$ws.Range("M257").Value = 845
$ws.Range("O257").Value = 'Región de O''Higgins'
$ws.Range("P257").Value = 845

$ws.Range("D258").Value = 44383
$ws.Range("J258").Value = 2000
$ws.Range("K258").Value = 800
$ws.Range("L258").Value = 800
$ws.Range("M258").Value = 800
$ws.Range("O258").Value = 'Región Metropolitana'
$ws.Range("P258").Value = 800

$ws.Range("D259").Value = 44307
$ws.Range("J259").Value = 850
$ws.Range("K259").Value = 1000
$ws.Range("L259").Value = 1000
$ws.Range("M259").Value = 1000
$ws.Range("O259").Value = 'Región Metropolitana'
$ws.Range("P259").Value = 1000

$ws.Range("D260").Value = 44307
$ws.Range("J260").Value = 1250
$ws.Range("K260").Value = 1000
$ws.Range("L260").Value = 1000
$ws.Range("M260").Value = 1000
$ws.Range("O260").Value = 'Región del Maule'
$ws.Range("P260").Value = 1000

$ws.Range("D261").Value = 44273
$ws.Range("J261").Value = 800
$ws.Range("K261").Value = 1000
$ws.Range("L261").Value = 1000
$ws.Range("M261").Value = 1000
$ws.Range("O261").Value = 'Provincia de Cautín'
$ws.Range("P261").Value = 1000

$ws.Range("D262").Value = 44273
$ws.Range("J262").Value = 1400
$ws.Range("K262").Value = 1000
$ws.Range("L262").Value = 1000
$ws.Range("M262").Value = 1000
$ws.Range("O262").Value = 'Región Metropolitana'
$ws.Range("P262").Value = 1000

$ws.Range("D263").Value = 44273
$ws.Range("J263").Value = 1200
$ws.Range("K263").Value = 1000
$ws.Range("L263").Value = 1000
$ws.Range("M263").Value = 1000
$ws.Range("O263").Value = 'Región del Maule'
$ws.Range("P263").Value = 1000

$ws.Range("D264").Value = 44433
$ws.Range("J264").Value = 1850
$ws.Range("K264").Value = 800
$ws.Range("L264").Value = 800
$ws.Range("M264").Value = 800
$ws.Range("O264").Value = 'Región de O''Higgins'
$ws.Range("P264").Value = 800

$ws.Range("D265").Value = 44302
$ws.Range("J265").Value = 500
$ws.Range("K265").Value = 800
$ws.Range("L265").Value = 900
$ws.Range("M265").Value = 860
$ws.Range("O265").Value = 'Región Metropolitana'
$ws.Range("P265").Value = 860

# Rows 266-268: new rows appended (266, 267) and the former last row shifted to 268
# Preserve the "Fecha" column's date number format (style index used by column D)
# on the newly-created rows, matching the existing rows above.
$dateFmt = $ws.Range("D253").NumberFormat
$ws.Range("D266").NumberFormat = $dateFmt
$ws.Range("D267").NumberFormat = $dateFmt
$ws.Range("D268").NumberFormat = $dateFmt

$ws.Range("A266").Value = 10
$ws.Range("B266").Value = 'Vega Modelo de Temuco'
$ws.Range("C266").Value = 'La Araucanía'
$ws.Range("D266").Value = 44302
$ws.Range("E266").Value = 9
$ws.Range("F266").Value = 100112008
$ws.Range("G266").Value = 'Coliflor'
$ws.Range("H266").Value = 'Sin especificar'
$ws.Range("I266").Value = 'Primera'
$ws.Range("J266").Value = 2700
$ws.Range("K266").Value = 800
$ws.Range("L266").Value = 900
$ws.Range("M266").Value = 844
$ws.Range("N266").Value = '$/unidad'
$ws.Range("O266").Value = 'Región del Maule'
$ws.Range("P266").Value = 844
$ws.Range("Q266").Value = 1
$ws.Range("R266").Value = 'Hortaliza'

$ws.Range("A267").Value = 10
$ws.Range("B267").Value = 'Vega Modelo de Temuco'
$ws.Range("C267").Value = 'La Araucanía'
$ws.Range("D267").Value = 44179
$ws.Range("E267").Value = 9
$ws.Range("F267").Value = 100112008
$ws.Range("G267").Value = 'Coliflor'
$ws.Range("H267").Value = 'Sin especificar'
$ws.Range("I267").Value = 'Primera'
$ws.Range("J267").Value = 2000
$ws.Range("K267").Value = 900
$ws.Range("L267").Value = 900
$ws.Range("M267").Value = 900
$ws.Range("N267").Value = '$/unidad'
$ws.Range("O267").Value = 'Región del Maule'
$ws.Range("P267").Value = 900
$ws.Range("Q267").Value = 1
$ws.Range("R267").Value = 'Hortaliza'

$ws.Range("A268").Value = 10
$ws.Range("B268").Value = 'Vega Modelo de Temuco'
$ws.Range("C268").Value = 'La Araucanía'
$ws.Range("D268").Value = 44491
$ws.Range("E268").Value = 9
$ws.Range("F268").Value = 100112008
$ws.Range("G268").Value = 'Coliflor'
$ws.Range("H268").Value = 'Sin especificar'
$ws.Range("I268").Value = 'Primera'
$ws.Range("J268").Value = 3200
$ws.Range("K268").Value = 800
$ws.Range("L268").Value = 900
$ws.Range("M268").Value = 839
$ws.Range("N268").Value = '$/unidad'
$ws.Range("O268").Value = 'Región Metropolitana'
$ws.Range("P268").Value = 839
$ws.Range("Q268").Value = 1
$ws.Range("R268").Value = 'Hortaliza'
